$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-28 Thursday" "2023-12-29 Friday"

Replace-Text "31×63=1953" "94×81=7614"
Replace-Text "54×95=5130" "46×40=1840"
Replace-Text "99×32=3168" "30×90=2700"
Replace-Text "13×42=546" "66×99=6534"
Replace-Text "51×32=1632" "31×35=1085"

Replace-Text "73×93=6789" "91×68=6188"
Replace-Text "25×42=1050" "68×90=6120"
Replace-Text "69×55=3795" "26×82=2132"
Replace-Text "84×58=4872" "61×78=4758"
Replace-Text "47×72=3384" "92×83=7636"

Replace-Text "29×12=348" "26×64=1664"
Replace-Text "60×24=1440" "12×27=324"
Replace-Text "43×33=1419" "70×34=2380"
Replace-Text "55×52=2860" "94×18=1692"
Replace-Text "12×17=204" "56×50=2800"

Replace-Text "31×54=1674" "88×48=4224"
Replace-Text "36×70=2520" "87×71=6177"
Replace-Text "15×88=1320" "88×37=3256"
Replace-Text "47×49=2303" "96×44=4224"
Replace-Text "77×65=5005" "55×64=3520"

Replace-Text "11×61=671" "60×57=3420"
Replace-Text "92×90=8280" "55×28=1540"
Replace-Text "59×89=5251" "31×66=2046"
Replace-Text "26×27=702" "48×84=4032"
Replace-Text "70×14=980" "96×49=4704"
